$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '56.701.87'
$ws.Range("E2").Value = '  +0.62%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.385.49'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '502.08'
$ws.Range("E5").Value = '  -1.66%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '131.19'
$ws.Range("E6").Value = '  +2.21%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.551'
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '2.392.22'
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0967'
$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("E11").Value = '  -1.17%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.321'
$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '4.55'
$ws.Range("E13").Value = '  -5.03%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '2.808.71'
$ws.Range("E14").Value = '  +1.05%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '56.655.90'
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '21.66'
$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '2.365.38'
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '10.16'
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '4.02'
$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '307.87'
$ws.Range("E21").Value = '  -1.10%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.25'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.86'
$ws.Range("E23").Value = '  +0.77%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '64.89'
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '0.373'
$ws.Range("E27").Value = '  -4.35%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.151'
$ws.Range("E28").Value = '  -0.85%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '7.39'
$ws.Range("E29").Value = '  +2.49%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '174.82'
$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.0₃0718'
$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.66'
$ws.Range("E32").Value = '  -1.06%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '5.89'
$ws.Range("E33").Value = '  -4.63%  '

$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '17.85'
$ws.Range("E37").Value = '  +1.24%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '1.19'
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("E39").Value = '  +2.68%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '36.83'
$ws.Range("E40").Value = '  +3.51%  '

$ws.Range("E41").Value = '  +2.51%  '

$ws.Range("E42").Value = '  +0.73%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '131.94'
$ws.Range("E43").Value = '  +4.17%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '3.35'
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '4.79'
$ws.Range("E45").Value = '  -1.57%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.562'
$ws.Range("E46").Value = '  -1.17%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.0905'
$ws.Range("E47").Value = '  +0.49%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '248.45'
$ws.Range("E48").Value = '  -2.62%  '

$ws.Range("E50").Value = '  +1.48%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '16.99'
$ws.Range("E51").Value = '  +8.48%  '
